$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "version 1 for no glob and p1, version 2 for p2 and p1_P2"
$ws.Range("B4").Select()
